# Apply the refreshed cryptos price/volume snapshot to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These new Price values parse as plain numbers; format as Text first so
# they are stored the same way as the rest of the (text) Price column.
$textCells = @("D5", "D6", "D11", "D12", "D14", "D19", "D20", "D21", "D22", "D24", "D25", "D28", "D30", "D32", "D34", "D36", "D37", "D40", "D42", "D43", "D48", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = '63.775.93'
$ws.Range("E2").Value = '  +0.74%  '

# Row 3 - Ethereum
$ws.Range("D3").Value = '3.092.51'
$ws.Range("E3").Value = '  -0.86%  '

# Row 4 - TetherUSD
$ws.Range("E4").Value = '  -0.03%  '

# Row 5 - BNB
$ws.Range("D5").Value = '540.68'
$ws.Range("E5").Value = '  -3.05%  '

# Row 6 - Solana
$ws.Range("D6").Value = '137.25'
$ws.Range("E6").Value = '  -1.56%  '

# Row 7 - USDC
$ws.Range("E7").Value = '  -0.05%  '

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = '3.086.30'
$ws.Range("E8").Value = '  -0.81%  '

# Row 9 - XRP
$ws.Range("E9").Value = '  -0.25%  '

# Row 10 - Dogecoin
$ws.Range("E10").Value = '  -2.90%  '

# Row 11 - Toncoin
$ws.Range("D11").Value = '6.33'
$ws.Range("E11").Value = '  -4.86%  '

# Row 12 - Cardano
$ws.Range("D12").Value = '0.459'
$ws.Range("E12").Value = '  -0.49%  '

# Row 13 - ShibaInu
$ws.Range("E13").Value = '  +3.93%  '

# Row 14 - Avalanche
$ws.Range("D14").Value = '34.87'
$ws.Range("E14").Value = '  -1.54%  '

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = '3.590.94'
$ws.Range("E15").Value = '  -1.13%  '

# Row 16 - WrappedBTC
$ws.Range("D16").Value = '63.780.73'
$ws.Range("E16").Value = '  +0.64%  '

# Row 17 - TRON
$ws.Range("E17").Value = '  +0.28%  '

# Row 18 - WrappedEther
$ws.Range("D18").Value = '3.091.58'
$ws.Range("E18").Value = '  -0.98%  '

# Row 19 - Polkadot
$ws.Range("D19").Value = '6.71'
$ws.Range("E19").Value = '  -0.60%  '

# Row 20 - BitcoinCash
$ws.Range("D20").Value = '489.71'
$ws.Range("E20").Value = '  -3.90%  '

# Row 21 - Chainlink
$ws.Range("D21").Value = '13.51'
$ws.Range("E21").Value = '  -1.01%  '

# Row 22 - Polygon
$ws.Range("D22").Value = '0.703'
$ws.Range("E22").Value = '  -1.44%  '

# Row 23 - Uniswap
$ws.Range("E23").Value = '  -1.89%  '

# Row 24 - Litecoin
$ws.Range("D24").Value = '79.90'
$ws.Range("E24").Value = '  +1.95%  '

# Row 25 - InternetComputer(DFINITY)
$ws.Range("D25").Value = '12.28'
$ws.Range("E25").Value = '  -1.67%  '

# Row 26 - Dai
$ws.Range("E26").Value = '  +0.14%  '

# Row 27 - PancakeSwap
$ws.Range("E27").Value = '  -1.72%  '

# Row 28 - RenderToken
$ws.Range("D28").Value = '8.33'
$ws.Range("E28").Value = '  -0.25%  '

# Row 29 - FirstDigitalUSD
$ws.Range("E29").Value = '  -0.13%  '

# Row 30 - EthereumClassic
$ws.Range("D30").Value = '26.32'
$ws.Range("E30").Value = '  -0.65%  '

# Row 31 - ImmutableX
$ws.Range("E31").Value = '  -2.92%  '

# Row 32 - Mantle
$ws.Range("D32").Value = '1.12'
$ws.Range("E32").Value = '  -0.38%  '

# Row 33 - Stacks
$ws.Range("E33").Value = '  -5.00%  '

# Row 34 - OKB
$ws.Range("D34").Value = '56.89'
$ws.Range("E34").Value = '  -4.36%  '

# Row 35 - NEARProtocol
$ws.Range("E35").Value = '  +4.59%  '

# Row 36 - now Filecoin (swapped with row 37)
$ws.Range("B36").Value = 'Filecoin'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D36").Value = '6.09'
$ws.Range("E36").Value = '  +1.29%  '

# Row 37 - now Bittensor (swapped with row 36)
$ws.Range("B37").Value = 'Bittensor'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D37").Value = '494.69'
$ws.Range("E37").Value = '  -7.26%  '

# Row 38 - Maker
$ws.Range("D38").Value = '3.304.74'
$ws.Range("E38").Value = '  +7.27%  '

# Row 39 - VeChain
$ws.Range("E39").Value = '  -4.13%  '

# Row 40 - Hedera
$ws.Range("D40").Value = '0.0802'
$ws.Range("E40").Value = '  +0.12%  '

# Row 41 - Kaspa
$ws.Range("E41").Value = '  -2.41%  '

# Row 42 - Cosmos
$ws.Range("D42").Value = '8.17'
$ws.Range("E42").Value = '  +0.22%  '

# Row 43 - dogwifhat
$ws.Range("D43").Value = '2.66'
$ws.Range("E43").Value = '  -3.47%  '

# Row 44 - TheGraph
$ws.Range("E44").Value = '  +0.70%  '

# Row 45 - USDe
$ws.Range("E45").Value = '  +0.04%  '

# Row 46 - Fetch.AI
$ws.Range("E46").Value = '  +0.42%  '

# Row 47 - PEPE
$ws.Range("E47").Value = '  +4.87%  '

# Row 48 - InjectiveProtocol
$ws.Range("D48").Value = '25.08'
$ws.Range("E48").Value = '  +2.52%  '

# Row 49 - Monero
$ws.Range("D49").Value = '121.94'
$ws.Range("E49").Value = '  -0.60%  '

# Row 50 - Stellar
$ws.Range("D50").Value = '0.110'
$ws.Range("E50").Value = '  +1.88%  '

# Row 51 - CoreDAO
$ws.Range("E51").Value = '  -4.47%  '
